$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.864.36'
$ws.Range("E2").Value = '  -0.96%  '

$ws.Range("D3").Value = '1.854.52'
$ws.Range("E3").Value = '  -0.58%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.43'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.94%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9998'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.07%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5046'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -2.15%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3650'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -2.89%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07172'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.16%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8895'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.50%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.66'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.06%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07508'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.75%  '

$ws.Range("D13").Value = '1.853.00'
$ws.Range("E13").Value = '  -0.77%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '91.60'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +2.70%  '

$ws.Range("E15").Value = '  -2.00%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.0000'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.00%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008522'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.42%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.05'
$ws.Range("D18").ClearFormats()

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9997'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.08%  '

$ws.Range("D20").Value = '26.897.36'
$ws.Range("E20").Value = '  -0.97%  '

$ws.Range("E21").Value = '  -0.20%  '

$ws.Range("D22").Value = '2.084.32'
$ws.Range("E22").Value = '  -1.79%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.31'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.74%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.440'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.57%  '

$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.33'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -3.24%  '

$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.797'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.65%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.81'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.19%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.053'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -5.15%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '112.68'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.26%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.643'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.21%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.655'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.85%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09209'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.17%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05093'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.17%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.984'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -3.85%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7353'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -2.72%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.147'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -2.35%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.227'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +6.84%  '

$ws.Range("E38").Value = '  -1.69%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.498'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.04%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.072'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.07%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5317'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.00%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '118.79'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +3.31%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.485'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.65%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.367'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.09%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1465'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.24%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4630'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.98%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9991'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.12%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.960'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.20%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.557'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.13%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.85'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.29%  '
